$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the decision table header text in A3 (expanded with additional clause)
$ws.Range("A3").Value = "DECISION TABLE FOR JUST ONE PLAYER'S MARBLES ON BOARD (not taking into affect other players playing at the same time) - so start here is only occuppied by one of playerX's 4 marbles, not another player yet…"

# Add two new documentation notes below the table (order of entry matters for shared-string indices)
$ws.Range("A18").Value = "next to take into consideration is other players"
$ws.Range("A16").Value = "can we jump ourselves? And what about shortcuts of the corners (i.e. land on an inside corner exactly & next time you roll 3, then can jump corner to corner)?"

# Update the selection to match the committed state (multi-area E6,E8,E10,E11 with E11 active)
$ws.Range("E6,E8,E10,E11").Select()
$ws.Range("E11").Activate()
